{"js": "// The dossier template used single-brace placeholders ({n_dossier},\n// {prenom}, {apogee}, {date}) for some fields while others had already\n// been migrated to double-brace ({{cne}}, {{email}}, {{nom_ar}}, \u2026).\n// This fixes the remaining single-brace placeholders so the generation\n// form (Jinja2-style double-brace templating) works for every field.\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"{n_dossier}\", replace: \"{{n_dossier}}\" },\n  { find: \"{prenom}\", replace: \"{{prenom}}\" },\n  { find: \"{apogee}\", replace: \"{{apogee}}\" },\n  { find: \"{date}\", replace: \"{{date}}\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The dossier template used single-brace placeholders ({n_dossier},\n# {prenom}, {apogee}, {date}) for some fields while others had already\n# been migrated to double-brace ({{cne}}, {{email}}, {{nom_ar}}, ...).\n# This fixes the remaining single-brace placeholders so the generation\n# form (Jinja2-style double-brace templating) works for every field.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$pairs = @(\n    @{ find = \"{n_dossier}\"; repl = \"{{n_dossier}}\" },\n    @{ find = \"{prenom}\";    repl = \"{{prenom}}\" },\n    @{ find = \"{apogee}\";    repl = \"{{apogee}}\" },\n    @{ find = \"{date}\";      repl = \"{{date}}\" }\n)\n\nforeach ($p in $pairs) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $p.find\n    $find.Replacement.Text = $p.repl\n    $find.Execute($p.find, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $p.repl, $wdReplaceAll) | Out-Null\n}\n"}
